$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 37726132
$ws.Range("I11").Value = 37726132
$ws.Range("K11").Value = 37726132
$ws.Range("M11").Value = -37725992

# Row 21
$ws.Range("H21").Value = 50017
$ws.Range("I21").Value = 50017
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 50017
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -49549
$ws.Range("N21").ClearContents()

# Row 23
$ws.Range("H23").Value = 50017
$ws.Range("I23").Value = 50017
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 50017
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -49783
$ws.Range("N23").ClearContents()

# Row 38
$ws.Range("H38").Value = 1142.5
$ws.Range("I38").Value = 173.71428
$ws.Range("J38").Value = 3403
$ws.Range("K38").Value = 521.14284
$ws.Range("L38").Value = 10209
$ws.Range("M38").Value = -149.14284
$ws.Range("N38").Value = -10953

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 18265.822
$ws.Range("I32").Value = 2089.1206
$ws.Range("J32").Value = 252828
$ws.Range("K32").Value = 2089.1206
$ws.Range("L32").Value = 252828
$ws.Range("M32").Value = -1802.1206
$ws.Range("N32").Value = -253402

# Row 45
$ws.Range("H45").Value = 954.5714
$ws.Range("I45").Value = 946.4
$ws.Range("J45").Value = 975
$ws.Range("K45").Value = 946.4
$ws.Range("L45").Value = 975
$ws.Range("M45").Value = -569.4
$ws.Range("N45").Value = -1729

# Row 74
$ws.Range("H74").Value = 4284.1465
$ws.Range("I74").Value = 1203.7715
$ws.Range("J74").Value = 22253
$ws.Range("K74").Value = 1203.7715
$ws.Range("L74").Value = 22253
$ws.Range("M74").Value = -329.7715000000001
$ws.Range("N74").Value = -24001

# Row 77
$ws.Range("H77").Value = 4284.1465
$ws.Range("I77").Value = 1203.7715
$ws.Range("J77").Value = 22253
$ws.Range("K77").Value = 6018.8575
$ws.Range("L77").Value = 111265
$ws.Range("M77").Value = -1650.8575
$ws.Range("N77").Value = -120001

# Row 122
$ws.Range("H122").Value = 1696
$ws.Range("I122").Value = 1681
$ws.Range("J122").Value = 1742.6666
$ws.Range("K122").Value = 5043
$ws.Range("L122").Value = 5227.9998
$ws.Range("M122").Value = -2593
$ws.Range("N122").Value = -10127.9998

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1323.3125
$ws.Range("I94").Value = 1137
$ws.Range("J94").Value = 1882.25
$ws.Range("K94").Value = 1137
$ws.Range("L94").Value = 1882.25
$ws.Range("M94").Value = -686
$ws.Range("N94").Value = -2784.25

$ws = $wb.Worksheets.Item("CRP")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 58
$ws.Range("H58").Value = 1996.1777
$ws.Range("I58").Value = 769.38464
$ws.Range("K58").Value = 769.38464
$ws.Range("M58").Value = -566.38464

# Row 64
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496

# Row 67
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716

# Row 132
$ws.Range("H132").Value = 2123.2678
$ws.Range("I132").Value = 1551.1777
$ws.Range("K132").Value = 4653.5331
$ws.Range("M132").Value = -2123.5331

# Row 134
$ws.Range("H134").Value = 2088.2456
$ws.Range("I134").Value = 1224.8223
$ws.Range("J134").Value = 5326.0835
$ws.Range("K134").Value = 3674.4669
$ws.Range("L134").Value = 15978.2505
$ws.Range("M134").Value = -1139.4669
$ws.Range("N134").Value = -21048.2505

# Row 136
$ws.Range("H136").Value = 1996.1777
$ws.Range("I136").Value = 769.38464
$ws.Range("K136").Value = 2308.15392
$ws.Range("M136").Value = 241.8460800000003

# Row 141
$ws.Range("H141").Value = 932500
$ws.Range("J141").Value = 1800000
$ws.Range("L141").Value = 1800000
$ws.Range("N141").Value = -1810360

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 24307.273
$ws.Range("I46").Value = 4980
$ws.Range("J46").Value = 26240
$ws.Range("K46").Value = 4980
$ws.Range("L46").Value = 26240
$ws.Range("M46").Value = -4824
$ws.Range("N46").Value = -26552

$ws = $wb.Worksheets.Item("LTW")
# Row 45
$ws.Range("H45").Value = 2050
$ws.Range("I45").Value = 2050
$ws.Range("K45").Value = 2050
$ws.Range("M45").Value = -1643

# Row 48
$ws.Range("H48").Value = 18333.334
$ws.Range("I48").Value = 18333.334
$ws.Range("K48").Value = 18333.334
$ws.Range("M48").Value = -17672.334

# Row 122
$ws.Range("H122").Value = 3075.4827
$ws.Range("I122").Value = 1917.6364
$ws.Range("J122").Value = 3783.0557
$ws.Range("K122").Value = 5752.9092
$ws.Range("L122").Value = 11349.1671
$ws.Range("M122").Value = -3302.9092
$ws.Range("N122").Value = -16249.1671

# Row 132
$ws.Range("H132").Value = 3692.3057
$ws.Range("I132").Value = 2517.68
$ws.Range("J132").Value = 6361.909
$ws.Range("K132").Value = 7553.039999999999
$ws.Range("L132").Value = 19085.727
$ws.Range("M132").Value = -5023.039999999999
$ws.Range("N132").Value = -24145.727

# Row 136
$ws.Range("H136").Value = 5403.731
$ws.Range("I136").Value = 3290.45
$ws.Range("J136").Value = 12448
$ws.Range("K136").Value = 9871.349999999999
$ws.Range("L136").Value = 37344
$ws.Range("M136").Value = -7321.349999999999
$ws.Range("N136").Value = -42444

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 12000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# Row 44
$ws.Range("H44").Value = 6000
$ws.Range("J44").Value = 6000
$ws.Range("L44").Value = 6000
$ws.Range("N44").Value = -7108

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 122
$ws.Range("H122").Value = 78629.46000000001
$ws.Range("J122").Value = 1938.125
$ws.Range("L122").Value = 5814.375
$ws.Range("N122").Value = -10714.375

# Row 126
$ws.Range("H126").Value = 112411.78
$ws.Range("I126").Value = 250525.25
$ws.Range("J126").Value = 1921
$ws.Range("K126").Value = 751575.75
$ws.Range("L126").Value = 5763
$ws.Range("M126").Value = -749105.75
$ws.Range("N126").Value = -10703

# Row 139
$ws.Range("H139").Value = 59866.668
$ws.Range("J139").Value = 59866.668
$ws.Range("L139").Value = 59866.668
$ws.Range("N139").Value = -70146.66800000001
